# Updates the cryptos listing sheet with refreshed price / volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.938.38"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.638.69"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.79"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.865.57"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "1.601.15"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.58"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "25.954.37"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.88"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  -1.70%  "
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.94"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  +2.66%  "
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("D37").Value = "1.138.07"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.52"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -0.95%  "
$ws.Range("E44").Value = "  -4.40%  "
$ws.Range("D45").Value = "1.774.37"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  +5.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.61"
$ws.Range("E47").Value = "  +1.60%  "
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.415"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  -0.72%  "

Write-Output "Applied 71 cell updates"
